# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---------------
# Copy the formatting of the existing last header cell (AC1, col 29)
# onto the three new header cells so they share the bold/border/centered
# style already used by the rest of row 1, then set their text.
$headerSrc = $ws.Cells.Item(1, 29)
$headerSrc.Copy()

$wins = $ws.Cells.Item(1, 30)
$wins.PasteSpecial(-4122)
$wins.Value = "Wins"

$losses = $ws.Cells.Item(1, 31)
$losses.PasteSpecial(-4122)
$losses.Value = "Losses"

$ties = $ws.Cells.Item(1, 32)
$ties.PasteSpecial(-4122)
$ties.Value = "Ties"

# --- Data rows (2-48): every player gets the team's season record -----
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
